$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3128873333333334
$ws.Range("H2").Value = 0.938662
$ws.Range("I2").Value = 0.1947636400758337
$ws.Range("J2").Value = 0.1947636400758337
$ws.Range("M2").Value = 25.69910333333333
$ws.Range("N2").Value = 77.09731
$ws.Range("O2").Value = 0.08761243344445813
$ws.Range("P2").Value = 0.08761243344445814
$ws.Range("Q2").Value = 8.040923911024445
$ws.Range("R2").Value = 72.36831519921999
$ws.Range("S2").Value = 0.01706371645354438
$ws.Range("T2").Value = 0.01706371645354438

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3128873333333334
$ws.Range("H3").Value = 0.938662
$ws.Range("I3").Value = 0.1947636400758337
$ws.Range("J3").Value = 0.1947636400758337
$ws.Range("O3").Value = 0.0005530844306649811
$ws.Range("P3").Value = 0.0005530844306649812
$ws.Range("Q3").Value = 0.05076117222755556
$ws.Range("R3").Value = 0.456850550048
$ws.Range("S3").Value = 0.0001077207369855818
$ws.Range("T3").Value = 0.0001077207369855818

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3128873333333334
$ws.Range("H4").Value = 0.938662
$ws.Range("I4").Value = 0.1947636400758337
$ws.Range("J4").Value = 0.1947636400758337
$ws.Range("M4").Value = 46.92720933333334
$ws.Range("N4").Value = 140.781628
$ws.Range("O4").Value = 0.1599825079935015
$ws.Range("P4").Value = 0.1599825079935016
$ws.Range("Q4").Value = 14.68292938908178
$ws.Range("R4").Value = 132.146364501736
$ws.Range("S4").Value = 0.03115877560527552
$ws.Range("T4").Value = 0.03115877560527552

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3128873333333334
$ws.Range("H5").Value = 0.938662
$ws.Range("I5").Value = 0.1947636400758337
$ws.Range("J5").Value = 0.1947636400758337
$ws.Range("M5").Value = 220.538579
$ws.Range("N5").Value = 661.615737
$ws.Range("O5").Value = 0.7518519741313753
$ws.Range("P5").Value = 0.7518519741313754
$ws.Range("Q5").Value = 69.00372788043266
$ws.Range("R5").Value = 621.0335509238939
$ws.Range("S5").Value = 0.1464334272800282
$ws.Range("T5").Value = 0.1464334272800282

$ws.Range("I6").Value = 0.2818555810746068
$ws.Range("J6").Value = 0.2818555810746068
$ws.Range("M6").Value = 25.69910333333333
$ws.Range("N6").Value = 77.09731
$ws.Range("O6").Value = 0.08761243344445813
$ws.Range("P6").Value = 0.08761243344445814
$ws.Range("Q6").Value = 11.63656255570111
$ws.Range("R6").Value = 104.72906300131
$ws.Range("S6").Value = 0.02469405333784806
$ws.Range("T6").Value = 0.02469405333784806

$ws.Range("I7").Value = 0.2818555810746068
$ws.Range("J7").Value = 0.2818555810746068
$ws.Range("O7").Value = 0.0005530844306649811
$ws.Range("P7").Value = 0.0005530844306649812
$ws.Range("S7").Value = 0.0001558899335883963
$ws.Range("T7").Value = 0.0001558899335883964

$ws.Range("I8").Value = 0.2818555810746068
$ws.Range("J8").Value = 0.2818555810746068
$ws.Range("M8").Value = 46.92720933333334
$ws.Range("N8").Value = 140.781628
$ws.Range("O8").Value = 0.1599825079935015
$ws.Range("P8").Value = 0.1599825079935016
$ws.Range("Q8").Value = 21.24865602853644
$ws.Range("R8").Value = 191.237904256828
$ws.Range("S8").Value = 0.0450919627522813
$ws.Range("T8").Value = 0.04509196275228131

$ws.Range("I9").Value = 0.2818555810746068
$ws.Range("J9").Value = 0.2818555810746068
$ws.Range("M9").Value = 220.538579
$ws.Range("N9").Value = 661.615737
$ws.Range("O9").Value = 0.7518519741313753
$ws.Range("P9").Value = 0.7518519741313754
$ws.Range("Q9").Value = 99.85994208405965
$ws.Range("R9").Value = 898.739478756537
$ws.Range("S9").Value = 0.211913675050889
$ws.Range("T9").Value = 0.211913675050889

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5847920000000001
$ws.Range("H10").Value = 1.754376
$ws.Range("I10").Value = 0.3640167129613011
$ws.Range("J10").Value = 0.3640167129613011
$ws.Range("M10").Value = 25.69910333333333
$ws.Range("N10").Value = 77.09731
$ws.Range("O10").Value = 0.08761243344445813
$ws.Range("P10").Value = 0.08761243344445814
$ws.Range("Q10").Value = 15.02863003650667
$ws.Range("R10").Value = 135.25767032856
$ws.Range("S10").Value = 0.03189239003699241
$ws.Range("T10").Value = 0.03189239003699242

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5847920000000001
$ws.Range("H11").Value = 1.754376
$ws.Range("I11").Value = 0.3640167129613011
$ws.Range("J11").Value = 0.3640167129613011
$ws.Range("O11").Value = 0.0005530844306649811
$ws.Range("P11").Value = 0.0005530844306649812
$ws.Range("Q11").Value = 0.09487353518933335
$ws.Range("R11").Value = 0.8538618167040001
$ws.Range("S11").Value = 0.0002013319764407391
$ws.Range("T11").Value = 0.0002013319764407391

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5847920000000001
$ws.Range("H12").Value = 1.754376
$ws.Range("I12").Value = 0.3640167129613011
$ws.Range("J12").Value = 0.3640167129613011
$ws.Range("M12").Value = 46.92720933333334
$ws.Range("N12").Value = 140.781628
$ws.Range("O12").Value = 0.1599825079935015
$ws.Range("P12").Value = 0.1599825079935016
$ws.Range("Q12").Value = 27.44265660045867
$ws.Range("R12").Value = 246.983909404128
$ws.Range("S12").Value = 0.05823630669109951
$ws.Range("T12").Value = 0.05823630669109952

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5847920000000001
$ws.Range("H13").Value = 1.754376
$ws.Range("I13").Value = 0.3640167129613011
$ws.Range("J13").Value = 0.3640167129613011
$ws.Range("M13").Value = 220.538579
$ws.Range("N13").Value = 661.615737
$ws.Range("O13").Value = 0.7518519741313753
$ws.Range("P13").Value = 0.7518519741313754
$ws.Range("Q13").Value = 128.969196690568
$ws.Range("R13").Value = 1160.722770215112
$ws.Range("S13").Value = 0.2736866842567684
$ws.Range("T13").Value = 0.2736866842567685

$ws.Range("G14").Value = 0.256018
$ws.Range("H14").Value = 0.768054
$ws.Range("I14").Value = 0.1593640658882584
$ws.Range("J14").Value = 0.1593640658882584
$ws.Range("M14").Value = 25.69910333333333
$ws.Range("N14").Value = 77.09731
$ws.Range("O14").Value = 0.08761243344445813
$ws.Range("P14").Value = 0.08761243344445814
$ws.Range("Q14").Value = 6.579433037193333
$ws.Range("R14").Value = 59.21489733473999
$ws.Range("S14").Value = 0.01396227361607328
$ws.Range("T14").Value = 0.01396227361607328

$ws.Range("G15").Value = 0.256018
$ws.Range("H15").Value = 0.768054
$ws.Range("I15").Value = 0.1593640658882584
$ws.Range("J15").Value = 0.1593640658882584
$ws.Range("O15").Value = 0.0005530844306649811
$ws.Range("P15").Value = 0.0005530844306649812
$ws.Range("Q15").Value = 0.04153499489066667
$ws.Range("R15").Value = 0.373814954016
$ws.Range("S15").Value = 0.00008814178365026392
$ws.Range("T15").Value = 0.00008814178365026392

$ws.Range("G16").Value = 0.256018
$ws.Range("H16").Value = 0.768054
$ws.Range("I16").Value = 0.1593640658882584
$ws.Range("J16").Value = 0.1593640658882584
$ws.Range("M16").Value = 46.92720933333334
$ws.Range("N16").Value = 140.781628
$ws.Range("O16").Value = 0.1599825079935015
$ws.Range("P16").Value = 0.1599825079935016
$ws.Range("Q16").Value = 12.01421027910134
$ws.Range("R16").Value = 108.127892511912
$ws.Range("S16").Value = 0.0254954629448452
$ws.Range("T16").Value = 0.0254954629448452

$ws.Range("G17").Value = 0.256018
$ws.Range("H17").Value = 0.768054
$ws.Range("I17").Value = 0.1593640658882584
$ws.Range("J17").Value = 0.1593640658882584
$ws.Range("M17").Value = 220.538579
$ws.Range("N17").Value = 661.615737
$ws.Range("O17").Value = 0.7518519741313753
$ws.Range("P17").Value = 0.7518519741313754
$ws.Range("Q17").Value = 56.46184591842201
$ws.Range("R17").Value = 508.156613265798
$ws.Range("S17").Value = 0.1198181875436896
$ws.Range("T17").Value = 0.1198181875436896

